$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.584.57"

$ws.Range("D3").Value = "1.925.69"
$ws.Range("E3").Value = "  +1.40%  "

$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = "0.9998"
$cell.Style = "Normal"
$ws.Range("E4").Value = "  +0.15%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "245.09"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +2.04%  "

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "0.9998"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  +0.12%  "

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.4733"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  -1.52%  "

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.2898"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  -2.14%  "

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.06821"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  +2.69%  "

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "106.52"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  +4.20%  "

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "18.44"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  -2.48%  "

$ws.Range("E12").Value = "  +1.32%  "

$ws.Range("D13").Value = "1.916.95"
$ws.Range("E13").Value = "  +1.24%  "

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "5.351"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  +3.55%  "

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "0.6689"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  +0.64%  "

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "292.69"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  -5.77%  "

$ws.Range("D17").Value = "30.607.95"
$ws.Range("E17").Value = "  -1.21%  "

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "0.000007639"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  +0.32%  "

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "13.02"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  -1.92%  "

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "0.9998"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  +0.13%  "

$ws.Range("D21").Value = "2.158.49"
$ws.Range("E21").Value = "  +1.76%  "

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "5.399"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  +3.65%  "

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "0.9996"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  +0.14%  "

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "6.260"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  -0.03%  "

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "9.369"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  -0.17%  "

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "168.39"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  +0.23%  "

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "21.21"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  +1.58%  "

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "2.121"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  +6.18%  "

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "0.1080"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  -4.47%  "

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "1.390"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  +2.64%  "

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "4.149"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  -1.51%  "

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "4.004"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  -0.39%  "

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "0.05064"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  -1.23%  "

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "0.7373"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  -2.03%  "

$ws.Range("E35").Value = "  -1.38%  "

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "0.02100"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  +5.36%  "

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "2.725"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  -0.39%  "

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "2.689"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  -1.03%  "

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "2.069"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  +0.12%  "

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "111.01"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  +1.90%  "

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "0.8753"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  -1.03%  "

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "5.921"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  +3.60%  "

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "0.4319"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  +1.54%  "

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "0.9997"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  +0.14%  "

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "68.04"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  -1.86%  "

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "49.27"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  +15.07%  "

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "7.239"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  -2.50%  "

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "9.305"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  -0.53%  "

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "0.1225"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  -1.13%  "

$ws.Range("E50").Value = "  +0.19%  "

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "0.2475"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  +8.98%  "
